# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-17 08:20:01
# Applies the attendance-data refresh to the B2 session analysis workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Column width: "Students" column (I) widened to fit the new "0/52" values.
# (13.1667 character-units round-trips through Excel's stored-width formula
# to a saved <col> width of exactly 14.)
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 13.166666666666666

# ---------------------------------------------------------------------------
# Recorded-By lists: the backup account now sorts first ("backup@backdoor.com"
# before "System"/"system"), and admin@admin.com now sorts after "System".
# ---------------------------------------------------------------------------
$ws.Range("G2").Value  = "backup@backdoor.com, System, system"
$ws.Range("G4").Value  = "backup@backdoor.com, System"
$ws.Range("G5").Value  = "backup@backdoor.com, System"
$ws.Range("G7").Value  = "System, admin@admin.com"
$ws.Range("G8").Value  = "backup@backdoor.com, System"
$ws.Range("G28").Value = "backup@backdoor.com, System, system"
$ws.Range("G30").Value = "backup@backdoor.com, System"
$ws.Range("G31").Value = "backup@backdoor.com, System"
$ws.Range("G33").Value = "System, admin@admin.com"
$ws.Range("G34").Value = "backup@backdoor.com, System"
$ws.Range("G54").Value = "backup@backdoor.com, System, system"
$ws.Range("G56").Value = "backup@backdoor.com, System"
$ws.Range("G57").Value = "backup@backdoor.com, System"
$ws.Range("G59").Value = "System, admin@admin.com"
$ws.Range("G60").Value = "backup@backdoor.com, System"
$ws.Range("G80").Value = "backup@backdoor.com, System"
$ws.Range("G81").Value = "backup@backdoor.com, System"
$ws.Range("G82").Value = "backup@backdoor.com, System"
$ws.Range("G106").Value = "backup@backdoor.com, System"
$ws.Range("G107").Value = "backup@backdoor.com, System"
$ws.Range("G108").Value = "backup@backdoor.com, System"
$ws.Range("G132").Value = "backup@backdoor.com, System"
$ws.Range("G133").Value = "backup@backdoor.com, System"
$ws.Range("G134").Value = "backup@backdoor.com, System"

# ---------------------------------------------------------------------------
# Recorded attendance counts that shifted for two sessions
# ---------------------------------------------------------------------------
$ws.Range("H50").Value = "7/57"
$ws.Range("H76").Value = "1/55"

# ---------------------------------------------------------------------------
# Class / group statistics recalculated after the data refresh.
# The percentages are stored as literal text (e.g. "99.4%"), not numbers, so
# a leading apostrophe is used to stop Excel from re-parsing them as numeric
# percentages.
# ---------------------------------------------------------------------------
$ws.Range("L6").Value  = 155
$ws.Range("L7").Value  = 1
$ws.Range("L9").Value  = "'99.4%"
$ws.Range("L10").Value = "'66.9%"

$ws.Range("O15").Value = 25
$ws.Range("P15").Value = 1
$ws.Range("R15").Value = "'96.2%"
$ws.Range("S15").Value = "'68.7%"
$ws.Range("S16").Value = "'62.3%"
$ws.Range("S17").Value = "'58.9%"

# ---------------------------------------------------------------------------
# Session 23 for group B2A (row 24) flips from Recorded to Not Recorded:
# the recorder is cleared, the attendance count drops to 0/52, the status
# text changes, and the row is highlighted pink (like the "Not Recorded"
# legend color) instead of the normal green "recorded" row color.
# ---------------------------------------------------------------------------
$ws.Range("G24").Value = ""
$ws.Range("H24").Value = "0/52"
$ws.Range("I24").Value = "Not Recorded"

$row24 = $ws.Range("A24:I24")
$row24.Interior.Color = 12695295   # RGB(255,182,193) pink, packed as BGR long
$row24.Font.Color = 0              # RGB(0,0,0) black
$row24.HorizontalAlignment = -4108 # xlCenter
$row24.VerticalAlignment = -4108   # xlCenter
